$d = $word.ActiveDocument

# 1) Merge the three runs describing the "thu ly" paragraph into a single
#    run, dropping the trailing period after ${noiDungToCao}.
$old = '${doiTuongGiaiQuyet} đã nhận đơn tố cáo của ${nguoiDungDon} đối với  …………………………………………………………………(2) về việc ${noiDungToCao}.'
$new = '${doiTuongGiaiQuyet} đã nhận đơn tố cáo của ${nguoiDungDon} đối với  …………………………………………………………………(2) về việc ${noiDungToCao}'
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# 2) Normal style: turn off "allow punctuation to overflow the text
#    boundary" (w:overflowPunct val false).
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.HangingPunctuation = $false
